# Regenerate sval data to filter save games: update computed stat columns
# (TB, d2S, K, IP, sum) for each row while leaving the Win flag column (F)
# and date labels (A) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 0.001754667048134761
    "C2" = 0.0001537489499301437
    "D2" = 16.98373111632243
    "E2" = 71517.89157740913
    "G2" = 71534.87721694144

    "B3" = 0.0003714022599530242
    "C3" = 0.004309184025731883
    "D3" = 3.082599426703578
    "E3" = 6.48142807727062
    "G3" = 9.568708090259884

    "B4" = 0.1554434735375247
    "C4" = 1766.335244827366
    "D4" = 0.7127328510149897
    "E4" = 71517.89157740913
    "G4" = 73285.09499856105

    "B5" = 0.7287194209349384
    "C5" = 1.65323645889881
    "D5" = 16.98373111632243
    "E5" = 6.48142807727062
    "G5" = 25.8471150734268
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
